$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 50. Excel will shift all the existing
# rows 50..107 down to 51..108 (and auto-extend the used range / dimension).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with a new weekly observation. It is
# identical to the record that used to sit in row 50 (now row 51) except for
# the date (column D) and the volume (column J).
$ws.Cells.Item(50, 1).Value = 4
$ws.Cells.Item(50, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(50, 3).Value = "Los Lagos"
$ws.Cells.Item(50, 4).Value = 44413
$ws.Cells.Item(50, 5).Value = 10
$ws.Cells.Item(50, 6).Value = 100112017
$ws.Cells.Item(50, 7).Value = "Apio"
$ws.Cells.Item(50, 8).Value = "Americana (o)"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 20
$ws.Cells.Item(50, 11).Value = 12000
$ws.Cells.Item(50, 12).Value = 12000
$ws.Cells.Item(50, 13).Value = 12000
$ws.Cells.Item(50, 14).Value = "$/docena de matas"
$ws.Cells.Item(50, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(50, 16).Value = 2000
$ws.Cells.Item(50, 17).Value = 6
$ws.Cells.Item(50, 18).Value = "Hortaliza"
